# Update stats for 2026-01 (row 26)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B26").Value = 6523
$ws.Range("C26").Value = 1013
$ws.Range("D26").Value = 6071701
$ws.Range("E26").Value = 930.8141959221217
$ws.Range("F26").Value = 10.24167652526618
$ws.Range("G26").Value = 7.537154989384298
$ws.Range("H26").Value = 26.44337787877456
